# "removed fiona code and finalized db_building"
#
# Semantic changes applied:
#  1. Master!A3  : "R1"   -> "R2"
#  2. Master!C3  : "COM2" -> "COM1"   (drops "COM2" from the shared-string
#                                      table since it becomes unused)
#  3. new_act2   : delete row 3 entirely (the leftover/duplicate data row)
#  4. new_act2!F2: "GLOBAL" -> "R1"
#  5. Active sheet/selection bookkeeping:
#       - Master becomes the active (selected) sheet/tab, with the
#         selection collapsed from G3:L3 to G3
#       - new_act1's selection moves from G2 to H27
#       - new_act2's selection moves from D15 to F2 (and it is no longer
#         the active tab)

$wb = $excel.ActiveWorkbook

$master  = $wb.Worksheets.Item("Master")
$act1    = $wb.Worksheets.Item("new_act1")
$act2    = $wb.Worksheets.Item("new_act2")

# --- Master sheet: update row 3 values -------------------------------
$master.Range("A3").Value = "R2"
$master.Range("C3").Value = "COM1"

# --- new_act2: drop the stray third row, fix F2 -----------------------
$act2.Range("F2").Value = "R1"
$act2.Rows.Item(3).Delete() | Out-Null

# --- Selections on each sheet ------------------------------------------
$act1.Range("H27").Select() | Out-Null
$act2.Range("F2").Select() | Out-Null

# Master selected/activated last so it ends up the active tab, with its
# own selection collapsed to G3 (was G3:L3).
$master.Range("G3").Select() | Out-Null
